$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "NLamps" column (M) ---

# Header cell M1: reuse the header formatting from L1 (border + style) so the
# existing style index is reused instead of creating a new one.
$ws.Range("M1").Value = "NLamps"
$ws.Range("L1").Copy()
$ws.Range("M1").PasteSpecial(-4122)

# Data values M2:M38 (number of lamps per configuration)
$values = @(1,2,1,2,1,2,1,2,3,4,1,2,3,4,1,2,3,4,1,2,3,4,2,4,6,8,2,4,8,5,11,7,8,5,5,3,2)
$arr = New-Object 'object[,]' $values.Length,1
for ($i = 0; $i -lt $values.Length; $i++) {
    $arr[$i,0] = $values[$i]
}
$ws.Range("M2:M38").Value = $arr

# Match the bordered/centered look used by the rest of the table (reuses an
# existing cell-format rather than defining a brand-new style).
$ws.Range("E2").Copy()
$ws.Range("M2:M38").PasteSpecial(-4122)

# Column M width ~ "9" characters, matching the other custom-width columns.
$ws.Columns("M").ColumnWidth = 8.166666666667

$excel.CutCopyMode = 0

# --- View / selection updates ---
$ws.Range("M4").Select() | Out-Null
